$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The records previously on rows 12/13/14 have been cyclically rotated:
#   new row 12 <- old row 13 data
#   new row 13 <- old row 14 data
#   new row 14 <- old row 12 data
# Only the columns that actually differ between the three source rows need
# to be rewritten (A, B, D, E, F, G, H, I, Q, R, S); the remaining columns
# were already identical across rows 12-14 and stay untouched.

# --- Row 12 (becomes old row 13's data) ---
$ws.Range("A12").Value = 111902037
$ws.Range("B12").Value = 90654
$ws.Range("E12").Value = 149
$ws.Range("F12").Value = "Tallgråticka"
$ws.Range("G12").Value = "Boletopsis grisea"
$ws.Range("H12").Value = "(Peck) Bondartsev & Singer"
$ws.Range("I12").Value = "'2"
$ws.Range("Q12").Value = 524868.6293626219
$ws.Range("R12").Value = 6867441.031870116

# --- Row 13 (becomes old row 14's data) ---
$ws.Range("A13").Value = 111902027
$ws.Range("B13").Value = 90660
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 4362
$ws.Range("F13").Value = "Blå taggsvamp"
$ws.Range("G13").Value = "Hydnellum caeruleum"
$ws.Range("H13").Value = "(Hornem.) P.Karst."
$ws.Range("I13").Value = "'5"
$ws.Range("Q13").Value = 524936.9216418237
$ws.Range("R13").Value = 6867321.952660743
$ws.Range("S13").Value = 25

# --- Row 14 (becomes old row 12's data) ---
$ws.Range("A14").Value = 111902029
$ws.Range("B14").Value = 88032
$ws.Range("D14").Value = "VU"
$ws.Range("E14").Value = 6276
$ws.Range("F14").Value = "Goliatmusseron"
$ws.Range("G14").Value = "Tricholoma matsutake"
$ws.Range("H14").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I14").Value = "'4"
$ws.Range("Q14").Value = 524971.6686743505
$ws.Range("R14").Value = 6867341.509407703
$ws.Range("S14").Value = 5
